$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.277.85"
$ws.Range("E2").Value = "  +5.46%  "
$ws.Range("D3").Value = "3.385.69"
$ws.Range("E3").Value = "  +6.21%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.56%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "3.391.66"
$ws.Range("E8").Value = "  +6.15%  "
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("E11").Value = "  +6.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.434"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").Value = "3.973.13"
$ws.Range("E13").Value = "  +6.18%  "
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000183"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.40%  "
$ws.Range("D17").Value = "63.418.69"
$ws.Range("E17").Value = "  +5.61%  "
$ws.Range("D18").Value = "3.412.75"
$ws.Range("E18").Value = "  +6.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.32%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.534"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.179"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000102"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +17.50%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  +7.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.57%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.63%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("E34").Value = "  +6.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0758"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.52%  "
$ws.Range("D41").Value = "2.901.90"
$ws.Range("E41").Value = "  +4.33%  "
$ws.Range("E42").Value = "  +3.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.761"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("E46").Value = "  +7.94%  "
$ws.Range("D47").Value = "3.437.39"
$ws.Range("E47").Value = "  +6.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "298.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +13.51%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.70%  "
